$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Bellarmine"
$ws.Range("B2").Value = "Atlantic Sun"
$ws.Range("C2").Value = "11-5"

$ws.Range("A3").Value = "Iowa State"
$ws.Range("B3").Value = "Big 12"
$ws.Range("C3").Value = "7-11"
